$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $r = $ws.Range($cellAddr)
    $r.Value = "`'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.402.18"
Set-TextValue "E2" "  +0.10%  "
Set-TextValue "D3" "1.845.89"
Set-TextValue "E3" "  -0.19%  "
Set-TextValue "D4" "0.9981"
Set-TextValue "E4" "  -0.20%  "
Set-TextValue "D5" "240.13"
Set-TextValue "E5" "  -0.19%  "
Set-TextValue "D6" "0.6292"
Set-TextValue "E6" "  +0.00%  "
Set-TextValue "D7" "0.9991"
Set-TextValue "E7" "  -0.06%  "
Set-TextValue "D8" "0.07466"
Set-TextValue "E8" "  -2.00%  "
Set-TextValue "D9" "0.2907"
Set-TextValue "D10" "24.52"
Set-TextValue "E10" "  -0.46%  "
Set-TextValue "D11" "0.07742"
Set-TextValue "D13" "5.002"
Set-TextValue "E13" "  -0.46%  "
Set-TextValue "D14" "0.6782"
Set-TextValue "D15" "0.00001048"
Set-TextValue "E15" "  -0.28%  "
Set-TextValue "D16" "82.13"
Set-TextValue "E16" "  -1.23%  "
Set-TextValue "D17" "6.193"
Set-TextValue "E17" "  +0.92%  "
Set-TextValue "D18" "29.394.00"
Set-TextValue "E18" "  +0.09%  "
Set-TextValue "D19" "228.55"
Set-TextValue "E19" "  -0.53%  "
Set-TextValue "D20" "12.35"
Set-TextValue "E20" "  +0.00%  "
Set-TextValue "E21" "  -0.05%  "
Set-TextValue "D22" "7.510"
Set-TextValue "E22" "  +0.43%  "
Set-TextValue "D23" "0.9987"
Set-TextValue "E23" "  -0.15%  "
Set-TextValue "D24" "159.25"
Set-TextValue "E24" "  +0.31%  "
Set-TextValue "D25" "8.509"
Set-TextValue "E25" "  +0.70%  "
Set-TextValue "D26" "0.1366"
Set-TextValue "E26" "  -1.93%  "
Set-TextValue "D27" "17.52"
Set-TextValue "D28" "0.06481"
Set-TextValue "E28" "  +15.34%  "
Set-TextValue "D29" "1.419"
Set-TextValue "E29" "  -2.47%  "
Set-TextValue "D30" "1.483"
Set-TextValue "E30" "  +0.57%  "
Set-TextValue "D31" "4.092"
Set-TextValue "E31" "  -0.55%  "
Set-TextValue "D32" "4.098"
Set-TextValue "E32" "  +1.07%  "
Set-TextValue "D33" "1.836"
Set-TextValue "E33" "  +0.06%  "
Set-TextValue "D34" "1.142"
Set-TextValue "E34" "  -1.47%  "
Set-TextValue "D35" "0.6996"
Set-TextValue "E35" "  +0.32%  "
Set-TextValue "D36" "2.584"
Set-TextValue "E36" "  -0.12%  "
Set-TextValue "D37" "1.263.78"
Set-TextValue "E37" "  +2.20%  "
Set-TextValue "D38" "0.01847"
Set-TextValue "E38" "  +1.92%  "
Set-TextValue "D39" "2.837"
Set-TextValue "E39" "  +4.00%  "
Set-TextValue "D40" "6.773"
Set-TextValue "E40" "  +5.51%  "
Set-TextValue "D41" "0.9291"
Set-TextValue "E41" "  +2.99%  "
Set-TextValue "D42" "0.9990"
Set-TextValue "E42" "  -0.12%  "
Set-TextValue "D43" "2.003.94"
Set-TextValue "E43" "  +1.12%  "
Set-TextValue "D44" "101.31"
Set-TextValue "E44" "  -0.26%  "
Set-TextValue "D45" "65.98"
Set-TextValue "E45" "  +0.45%  "
Set-TextValue "D46" "1.735"
Set-TextValue "E46" "  +2.89%  "
Set-TextValue "D47" "7.084"
Set-TextValue "E47" "  -1.27%  "
Set-TextValue "E48" "  +0.06%  "
Set-TextValue "B49" "BabyDogeCoin"
Set-TextValue "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.00000000116"
Set-TextValue "E49" "  +1.55%  "
Set-TextValue "B50" "EnergySwap"
Set-TextValue "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "9.032"
Set-TextValue "E50" "  +0.03%  "
Set-TextValue "D51" "0.3946"
Set-TextValue "E51" "  -1.32%  "
